$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2-86 from 45203 to 45204
for ($r = 2; $r -le 86; $r++) {
    $ws.Cells.Item($r, 3).Value = 45204
}

# Row 86 gains an explicit row height (15, custom height)
$ws.Rows.Item(86).RowHeight = 15

# Append new row 87 with the new logging notice entry
$ws.Range("A87").Value = "A 47646-2023"

$ws.Range("B87").Value = 45203
$ws.Range("B87").NumberFormat = "YYYY-MM-DD"

$ws.Range("C87").Value = 45204
$ws.Range("C87").NumberFormat = "YYYY-MM-DD"

$ws.Range("D87").Value = "VÄRMLANDS LÄN"
$ws.Range("E87").Value = "MUNKFORS"

$ws.Range("G87").Value = 1.9
$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = 0
$ws.Range("N87").Value = 0
$ws.Range("O87").Value = 0
$ws.Range("P87").Value = 0
$ws.Range("Q87").Value = 0

$ws.Range("R87").Value = ""
$ws.Range("R87").WrapText = $true
